# Raitha Dinachari.xlsx - "Add files via upload" commit replay
#
# Summary of the edit (per the OOXML diff):
#  - "Debt" tracker rows 20 & 21 on "Daily Expenditure" (previously blank)
#    get filled in with two new loans ("Bava" and "Harish Land Loan").
#  - Four new daily-expenditure rows (44-47) are appended, three of which
#    carry real transactions ("Ginger" / "Sales" / "Buy" / "Medicine") and
#    the last is a placeholder row with only the SlNo/Date/Month/Type filled.
#  - Two brand-new shared strings ("Buy", "Bava") are introduced by the
#    above edits; Excel appends these to sharedStrings.xml automatically
#    when we write them as cell text.
#  - The saved view state changed: "Daily Pivot" scrolled to C1 and
#    "Daily Expenditure" scrolled to M5 with the active cell at P33.
#  - Every TODAY()-driven formula (interest/ageing calculations, SUMIFs,
#    etc.) simply recalculates against the new "current date" - no manual
#    value pokes are required for those, the engine's automatic recalc
#    (xlAutomatic, matches run_com's post-script recalc) takes care of it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily Expenditure")

# --- Debt tracker: row 20 ("Bava") ---------------------------------------
$ws.Range("J20").Value = 44645
$ws.Range("K20").Value = "Debt"
$ws.Range("L20").Value = "Bava"
$ws.Range("M20").Value = 150000
$ws.Range("N20").Value = 1

# --- Debt tracker: row 21 ("Harish Land Loan") ---------------------------
$ws.Range("J21").Value = 44651
$ws.Range("K21").Value = "Debt"
$ws.Range("L21").Value = "Harish Land Loan"
$ws.Range("M21").Value = 750000
$ws.Range("N21").Value = 1

# --- Daily expenditure log: finish row 44, append rows 45-47 -------------
# Row 44 already exists (SlNo/Date blank, Expense type + TEXT() formula in
# place) - just fill in the transaction that was recorded for it.
$ws.Range("A44").Value = 44
$ws.Range("B44").Value = 44643
$ws.Range("E44").Value = "Ginger"
$ws.Range("F44").Value = "Sales"
$ws.Range("G44").Value = 11900

# New rows: copy row 44's number formats/styles down first so the new
# cells pick up the same styling (s="2"/"20" etc.) the workbook already
# uses for this table, then fill in the values.
$ws.Range("A44:G44").Copy()
$ws.Range("A45:G47").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A45").Value = 45
$ws.Range("B45").Value = 44643
$ws.Range("C45").Formula = '=TEXT(B45,"mmm")'
$ws.Range("D45").Value = "Expense"
$ws.Range("E45").Value = "Ginger"
$ws.Range("F45").Value = "Buy"
$ws.Range("G45").Value = 20000

$ws.Range("A46").Value = 45
$ws.Range("B46").Value = 44650
$ws.Range("C46").Formula = '=TEXT(B46,"mmm")'
$ws.Range("D46").Value = "Expense"
$ws.Range("E46").Value = "Ginger"
$ws.Range("F46").Value = "Medicine"
$ws.Range("G46").Value = 900

$ws.Range("A47").Value = 45
$ws.Range("B47").Value = 44650
$ws.Range("C47").Formula = '=TEXT(B47,"mmm")'
$ws.Range("D47").Value = "Expense"

# --- Restore saved view/selection state -----------------------------------
$wsPivot = $wb.Worksheets.Item("Daily Pivot")
$wsPivot.Range("C1").Select()

$ws.Range("M5").Select()
$ws.Range("P33").Select()
